$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2900.628494009878
$ws.Range("E2").Value = 290490.7128553874
$ws.Range("G2").Value = 80959.25712662051
$ws.Range("I2").Value = 149451.0834652955
$ws.Range("L2").Value = 509988.6069102
$ws.Range("M2").Value = 112287.0813999
$ws.Range("N2").Value = 71616.34392528504
$ws.Range("O2").Value = 66836.36011669762

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 2196.160489230463
$ws.Range("B2").Value = 35136.73892605074
$ws.Range("E2").Value = 164484.8714791987
$ws.Range("I2").Value = 161693.2696809839
$ws.Range("L2").Value = 94361.20555763146
$ws.Range("M2").Value = 61433.01601085002
$ws.Range("N2").Value = 19369.86334567976
$ws.Range("O2").Value = 11611.14661559149

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 19330.26082556661
$ws.Range("B2").Value = 19240.78555357244
$ws.Range("E2").Value = 120007.87107624
$ws.Range("I2").Value = 187102.6412512376
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 60144.33064563009
$ws.Range("N2").Value = 44280.52571750963
$ws.Range("O2").Value = 44977.21926178802

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("O2").Value = 19878.02997089422
